$d = $word.ActiveDocument

# The licence paragraph reads:
#   "Except where otherwise noted, this work by The Shodor Education
#    Foundation, Inc. is licensed under CC BY-NC 4.0. To view a copy of
#    this license, visit "
#   <break>
#   https://creativecommons.org/licenses/by-nc/4.0   (hyperlink)
#
# The project switched its licence from CC BY-NC 4.0 to CC BY-SA 4.0, so
# both the visible "CC BY-NC" label and the creativecommons.org hyperlink
# (its address and its display text) need to change to "CC BY-SA" / "by-sa".

# 1. Update the visible licence name in the body text.
$d.Content.Find.Execute("CC BY-NC 4.0", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "CC BY-SA 4.0", 2) | Out-Null

# 2. Update the Creative Commons hyperlink: both where it points and what
#    it displays, so the two stay in sync.
foreach ($h in $d.Hyperlinks) {
    if ($h.Address -eq "https://creativecommons.org/licenses/by-nc/4.0") {
        $h.TextToDisplay = "https://creativecommons.org/licenses/by-sa/4.0"
        $h.Address = "https://creativecommons.org/licenses/by-sa/4.0"
    }
}
